# Showcase App.docx - merge "dev" edits:
#  - drop now-unneeded w:proofErr spell-check bookmarks around a few runs
#    (and coalesce the runs they used to split into a single run each)
#  - "Login" -> "Login / Logout" (as a distinct appended run)
#  - "Profile" -> "Profile" + "s" (as a distinct appended run, i.e. "Profiles")
#  - replace the "Friends / Followers List" bullet with three bullets:
#    "Friends List", "Search", "Post Creation" (inserted right after Profile(s))
#  - remove the old "Search" bullet that used to sit right after "Post Feed"
#  - remove the trailing empty bullet paragraph after "Light / Dark mode"

$d = $word.ActiveDocument

# Replace a whole paragraph's contents with a single clean run (no proofErr,
# no leftover run-splits) while preserving the paragraph's own formatting
# (pStyle/numPr/etc). We do this by splitting a fresh, identically-formatted
# paragraph off of the front of the target paragraph, filling it with the
# desired text, and then deleting the original (now-redundant) paragraph,
# pilcrow included, so none of its old runs/proofErr markers survive.
function Set-CleanParagraphText($doc, $paraIndex, $newText) {
    $p = $doc.Paragraphs.Item($paraIndex)
    $start = $p.Range.Start
    $ins = $doc.Range($start, $start)
    $ins.InsertBefore("`r")

    $newPara = $doc.Paragraphs.Item($paraIndex)
    $newPara.Range.Text = $newText

    $oldPara = $doc.Paragraphs.Item($paraIndex + 1)
    $oldPara.Range.Delete()
}

# Append `newText` to a paragraph as a brand-new w:r (rather than merging it
# into the existing trailing run). We do this by splitting the text onto a
# temporary new paragraph right after the target one, then deleting the
# pilcrow between them so the two merge back into a single paragraph that
# keeps both runs distinct.
function Add-RunAfterParagraph($doc, $paraIndex, $newText) {
    $p = $doc.Paragraphs.Item($paraIndex)
    $endPos = $p.Range.End - 1
    $insPoint = $doc.Range($endPos, $endPos)
    $insPoint.InsertAfter("`r" + $newText)

    $p2 = $doc.Paragraphs.Item($paraIndex)
    $pilcrowPos = $p2.Range.End - 1
    $pilcrowRange = $doc.Range($pilcrowPos, $pilcrowPos + 1)
    $pilcrowRange.Delete()
}

# --- 1. Drop proofErr wrapping / coalesce runs (paragraph count unaffected) ---
Set-CleanParagraphText $d 5 "ExpressJS"
Set-CleanParagraphText $d 11 "Formik + Yup for form and validation"
Set-CleanParagraphText $d 14 "React Dropzone for image upload"
Set-CleanParagraphText $d 17 "ExpressJS for framework"
Set-CleanParagraphText $d 20 "Multer for file upload"

# --- 2. "Login" -> "Login" + " / Logout" ; "Profile" -> "Profile" + "s" ---
Add-RunAfterParagraph $d 23 " / Logout"
Add-RunAfterParagraph $d 24 "s"

# --- 3. Replace "Friends / Followers List" with 3 new bullets ---
$friendsPara = $d.Paragraphs.Item(25)
$start = $friendsPara.Range.Start
$ins = $d.Range($start, $start)
$ins.InsertBefore("Friends List`rSearch`rPost Creation`r")
# the original "Friends / Followers List" paragraph got pushed down by 3
$oldFriendsPara = $d.Paragraphs.Item(28)
$oldFriendsPara.Range.Delete()

# --- 4. Remove the old "Search" bullet that followed "Post Feed" ---
# After step 3, order is: ... Post Feed(28) Search(29) Post Likes(30) ...
$oldSearchPara = $d.Paragraphs.Item(29)
$oldSearchPara.Range.Delete()

# --- 5. Remove the empty bullet paragraph after "Light / Dark mode" ---
# After step 4, order is: ... Light / Dark mode(32) <empty>(33) <drawing>(34) ...
$emptyPara = $d.Paragraphs.Item(33)
$emptyPara.Range.Delete()
